# Update "想去人数" (F column) figures that changed between crawl runs.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 1786
$wsExhibit.Range("F9").Value = 1732
$wsExhibit.Range("F16").Value = 12801
$wsExhibit.Range("F27").Value = 29

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1786
$wsAll.Range("F14").Value = 1732
$wsAll.Range("F22").Value = 12801
$wsAll.Range("F25").Value = 10
$wsAll.Range("F37").Value = 29
